$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update password column (B) for all data rows to the new (typo'd) value
$ws.Range("B2:B7").Value = "secret_sauc"

# Update exp column (C) for all data rows to "Fail"
$ws.Range("C2:C7").Value = "Fail"

# Reflect the new active cell selection recorded in the saved file
$ws.Range("G6").Select()
